$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 1000  # H12: 1200 -> 1000
$ws.Cells.Item(12, 9).Value = 1000  # I12: 1200 -> 1000
$ws.Cells.Item(12, 11).Value = 1000  # K12: 1200 -> 1000
$ws.Cells.Item(12, 13).Value = -830  # M12: -1030 -> -830
$ws.Cells.Item(17, 8).Value = 625949.56  # H17: 667612.75 -> 625949.5600000001
$ws.Cells.Item(17, 10).Value = 1112021.5  # J17: 1250898.9 -> 1112021.5
$ws.Cells.Item(17, 12).Value = 3336064.5  # L17: 3752696.7 -> 3336064.5
$ws.Cells.Item(17, 14).Value = -3336400.5  # N17: -3753032.7 -> -3336400.5
$ws.Cells.Item(53, 8).Value = 2754.8696  # H53: 2640.125 -> 2754.8696
$ws.Cells.Item(53, 9).Value = 96  # I53: 88.083336 -> 96
$ws.Cells.Item(53, 11).Value = 96  # K53: 88.083336 -> 96
$ws.Cells.Item(53, 13).Value = 541  # M53: 548.916664 -> 541
$ws.Cells.Item(76, 8).Value = 18551492  # H76: 16696651 -> 18551492
$ws.Cells.Item(76, 9).Value = 3433  # I76: 3295.4 -> 3433
$ws.Cells.Item(76, 11).Value = 3433  # K76: 3295.4 -> 3433
$ws.Cells.Item(76, 13).Value = -3118  # M76: -2980.4 -> -3118
$ws.Cells.Item(79, 8).Value = 18551492  # H79: 16696651 -> 18551492
$ws.Cells.Item(79, 9).Value = 3433  # I79: 3295.4 -> 3433
$ws.Cells.Item(79, 11).Value = 3433  # K79: 3295.4 -> 3433
$ws.Cells.Item(79, 13).Value = -2341  # M79: -2203.4 -> -2341
$ws.Cells.Item(86, 8).Value = 8721880  # H86: 8358551.5 -> 8721880
$ws.Cells.Item(86, 9).Value = 5141.6665  # I86: 4827.5 -> 5141.6665
$ws.Cells.Item(86, 11).Value = 5141.6665  # K86: 4827.5 -> 5141.6665
$ws.Cells.Item(86, 13).Value = -4018.6665  # M86: -3704.5 -> -4018.6665
$ws.Cells.Item(89, 8).Value = 8721880  # H89: 8358551.5 -> 8721880
$ws.Cells.Item(89, 9).Value = 5141.6665  # I89: 4827.5 -> 5141.6665
$ws.Cells.Item(89, 11).Value = 25708.3325  # K89: 24137.5 -> 25708.3325
$ws.Cells.Item(89, 13).Value = -20092.3325  # M89: -18521.5 -> -20092.3325
$ws.Cells.Item(98, 8).Value = 2044.2413  # H98: 2047.8966 -> 2044.2413
$ws.Cells.Item(98, 10).Value = 1428.4286  # J98: 1443.5714 -> 1428.4286
$ws.Cells.Item(98, 12).Value = 1428.4286  # L98: 1443.5714 -> 1428.4286
$ws.Cells.Item(98, 14).Value = -4424.4286  # N98: -4439.5714 -> -4424.4286
$ws.Cells.Item(100, 8).Value = 3430.9412  # H100: 3708.0667 -> 3430.9412
$ws.Cells.Item(100, 9).Value = 3573.1428  # I100: 3943.25 -> 3573.1428
$ws.Cells.Item(100, 11).Value = 3573.1428  # K100: 3943.25 -> 3573.1428
$ws.Cells.Item(100, 13).Value = -3032.1428  # M100: -3402.25 -> -3032.1428
$ws.Cells.Item(122, 8).Value = 2044.2413  # H122: 2047.8966 -> 2044.2413
$ws.Cells.Item(122, 10).Value = 1428.4286  # J122: 1443.5714 -> 1428.4286
$ws.Cells.Item(122, 12).Value = 4285.2858  # L122: 4330.7142 -> 4285.2858
$ws.Cells.Item(122, 14).Value = -9185.2858  # N122: -9230.7142 -> -9185.2858
$ws.Cells.Item(138, 8).Value = 3176.691  # H138: 3143.8728 -> 3176.691
$ws.Cells.Item(138, 9).Value = 2274.8696  # I138: 2238.2083 -> 2274.8696
$ws.Cells.Item(138, 10).Value = 3824.875  # J138: 3845.0322 -> 3824.875
$ws.Cells.Item(138, 11).Value = 6824.6088  # K138: 6714.624899999999 -> 6824.6088
$ws.Cells.Item(138, 12).Value = 11474.625  # L138: 11535.0966 -> 11474.625
$ws.Cells.Item(138, 13).Value = -1684.6088  # M138: -1574.624899999999 -> -1684.6088
$ws.Cells.Item(138, 14).Value = -21754.625  # N138: -21815.0966 -> -21754.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 74500  # H44: 75000 -> 74500
$ws.Cells.Item(44, 10).Value = 74500  # J44: 75000 -> 74500
$ws.Cells.Item(44, 12).Value = 74500  # L44: 75000 -> 74500
$ws.Cells.Item(44, 14).Value = -75476  # N44: -75976 -> -75476
$ws.Cells.Item(61, 8).Value = 2001368.6  # H61: 2501648 -> 2001368.6
$ws.Cells.Item(61, 9).Value = 1251335.9  # I61: 2001437.2 -> 1251335.9
$ws.Cells.Item(61, 10).Value = 5001499.5  # J61: 3335332.2 -> 5001499.5
$ws.Cells.Item(61, 11).Value = 1251335.9  # K61: 2001437.2 -> 1251335.9
$ws.Cells.Item(61, 12).Value = 5001499.5  # L61: 3335332.2 -> 5001499.5
$ws.Cells.Item(61, 13).Value = -1251123.9  # M61: -2001225.2 -> -1251123.9
$ws.Cells.Item(61, 14).Value = -5001923.5  # N61: -3335756.2 -> -5001923.5
$ws.Cells.Item(110, 8).Value = 83334850  # H110: 76924560 -> 83334850
$ws.Cells.Item(110, 9).Value = 90910560  # I110: 83334770 -> 90910560
$ws.Cells.Item(110, 11).Value = 90910560  # K110: 83334770 -> 90910560
$ws.Cells.Item(110, 13).Value = -90908515  # M110: -83332725 -> -90908515
$ws.Cells.Item(122, 8).Value = 6234.1665  # H122: 5607.5 -> 6234.1665
$ws.Cells.Item(122, 9).Value = 2889.75  # I122: 2542.3333 -> 2889.75
$ws.Cells.Item(122, 11).Value = 8669.25  # K122: 7626.999899999999 -> 8669.25
$ws.Cells.Item(122, 13).Value = -6219.25  # M122: -5176.999899999999 -> -6219.25
$ws.Cells.Item(126, 8).Value = 10040  # H126: 15289.857 -> 10040
$ws.Cells.Item(126, 9).Value = 10040  # I126: 15289.857 -> 10040
$ws.Cells.Item(126, 11).Value = 30120  # K126: 45869.571 -> 30120
$ws.Cells.Item(126, 13).Value = -27650  # M126: -43399.571 -> -27650
$ws.Cells.Item(136, 8).Value = 2001368.6  # H136: 2501648 -> 2001368.6
$ws.Cells.Item(136, 9).Value = 1251335.9  # I136: 2001437.2 -> 1251335.9
$ws.Cells.Item(136, 10).Value = 5001499.5  # J136: 3335332.2 -> 5001499.5
$ws.Cells.Item(136, 11).Value = 3754007.7  # K136: 6004311.6 -> 3754007.7
$ws.Cells.Item(136, 12).Value = 15004498.5  # L136: 10005996.6 -> 15004498.5
$ws.Cells.Item(136, 13).Value = -3751457.7  # M136: -6001761.6 -> -3751457.7
$ws.Cells.Item(136, 14).Value = -15009598.5  # N136: -10011096.6 -> -15009598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 587.5  # H22: 225 -> 587.5
$ws.Cells.Item(22, 9).Value = 650  # I22: 225 -> 650
$ws.Cells.Item(22, 10).Value = 400  # J22: 0 -> 400
$ws.Cells.Item(22, 11).Value = 650  # K22: 225 -> 650
$ws.Cells.Item(22, 12).Value = 400  # L22: 0 -> 400
$ws.Cells.Item(22, 13).Value = -477  # M22: -52 -> -477
$ws.Cells.Item(22, 14).Value = -746  # N22: None -> -746
$ws.Cells.Item(134, 8).Value = 3021  # H134: 2983 -> 3021
$ws.Cells.Item(134, 9).Value = 2381.8333  # I134: 2475 -> 2381.8333
$ws.Cells.Item(134, 10).Value = 4299.3335  # J134: 3999 -> 4299.3335
$ws.Cells.Item(134, 11).Value = 7145.499899999999  # K134: 7425 -> 7145.499899999999
$ws.Cells.Item(134, 12).Value = 12898.0005  # L134: 11997 -> 12898.0005
$ws.Cells.Item(134, 13).Value = -4610.499899999999  # M134: -4890 -> -4610.499899999999
$ws.Cells.Item(134, 14).Value = -17968.0005  # N134: -17067 -> -17968.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 46709.062  # H74: 48713.547 -> 46709.062
$ws.Cells.Item(74, 9).Value = 20571.25  # I74: 20761.666 -> 20571.25
$ws.Cells.Item(74, 10).Value = 55421.668  # J74: 59195.5 -> 55421.668
$ws.Cells.Item(74, 11).Value = 20571.25  # K74: 20761.666 -> 20571.25
$ws.Cells.Item(74, 12).Value = 55421.668  # L74: 59195.5 -> 55421.668
$ws.Cells.Item(74, 13).Value = -19697.25  # M74: -19887.666 -> -19697.25
$ws.Cells.Item(74, 14).Value = -57169.668  # N74: -60943.5 -> -57169.668
$ws.Cells.Item(77, 8).Value = 46709.062  # H77: 48713.547 -> 46709.062
$ws.Cells.Item(77, 9).Value = 20571.25  # I77: 20761.666 -> 20571.25
$ws.Cells.Item(77, 10).Value = 55421.668  # J77: 59195.5 -> 55421.668
$ws.Cells.Item(77, 11).Value = 61713.75  # K77: 62284.99800000001 -> 61713.75
$ws.Cells.Item(77, 12).Value = 166265.004  # L77: 177586.5 -> 166265.004
$ws.Cells.Item(77, 13).Value = -57345.75  # M77: -57916.99800000001 -> -57345.75
$ws.Cells.Item(77, 14).Value = -175001.004  # N77: -186322.5 -> -175001.004
$ws.Cells.Item(86, 8).Value = 5732.125  # H86: 4925.048 -> 5732.125
$ws.Cells.Item(86, 9).Value = 5513.0835  # I86: 5265.923 -> 5513.0835
$ws.Cells.Item(86, 10).Value = 6389.25  # J86: 4371.125 -> 6389.25
$ws.Cells.Item(86, 11).Value = 5513.0835  # K86: 5265.923 -> 5513.0835
$ws.Cells.Item(86, 12).Value = 6389.25  # L86: 4371.125 -> 6389.25
$ws.Cells.Item(86, 13).Value = -4390.0835  # M86: -4142.923 -> -4390.0835
$ws.Cells.Item(86, 14).Value = -8635.25  # N86: -6617.125 -> -8635.25
$ws.Cells.Item(89, 8).Value = 5732.125  # H89: 4925.048 -> 5732.125
$ws.Cells.Item(89, 9).Value = 5513.0835  # I89: 5265.923 -> 5513.0835
$ws.Cells.Item(89, 10).Value = 6389.25  # J89: 4371.125 -> 6389.25
$ws.Cells.Item(89, 11).Value = 27565.4175  # K89: 26329.615 -> 27565.4175
$ws.Cells.Item(89, 12).Value = 31946.25  # L89: 21855.625 -> 31946.25
$ws.Cells.Item(89, 13).Value = -21949.4175  # M89: -20713.615 -> -21949.4175
$ws.Cells.Item(89, 14).Value = -43178.25  # N89: -33087.625 -> -43178.25
$ws.Cells.Item(99, 8).Value = 2515  # H99: 2149.1177 -> 2515
$ws.Cells.Item(99, 9).Value = 2333.75  # I99: 1980.8334 -> 2333.75
$ws.Cells.Item(99, 10).Value = 2998.3333  # J99: 2553 -> 2998.3333
$ws.Cells.Item(99, 11).Value = 2333.75  # K99: 1980.8334 -> 2333.75
$ws.Cells.Item(99, 12).Value = 2998.3333  # L99: 2553 -> 2998.3333
$ws.Cells.Item(99, 13).Value = -835.75  # M99: -482.8334 -> -835.75
$ws.Cells.Item(99, 14).Value = -5994.3333  # N99: -5549 -> -5994.3333
$ws.Cells.Item(126, 8).Value = 2515  # H126: 2149.1177 -> 2515
$ws.Cells.Item(126, 9).Value = 2333.75  # I126: 1980.8334 -> 2333.75
$ws.Cells.Item(126, 10).Value = 2998.3333  # J126: 2553 -> 2998.3333
$ws.Cells.Item(126, 11).Value = 7001.25  # K126: 5942.5002 -> 7001.25
$ws.Cells.Item(126, 12).Value = 8994.999899999999  # L126: 7659 -> 8994.999899999999
$ws.Cells.Item(126, 13).Value = -4531.25  # M126: -3472.5002 -> -4531.25
$ws.Cells.Item(126, 14).Value = -13934.9999  # N126: -12599 -> -13934.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 375  # H7: 425 -> 375
$ws.Cells.Item(7, 9).Value = 375  # I7: 425 -> 375
$ws.Cells.Item(7, 11).Value = 1125  # K7: 1275 -> 1125
$ws.Cells.Item(7, 13).Value = -1013  # M7: -1163 -> -1013
$ws.Cells.Item(80, 8).Value = 2799.6667  # H80: 3469 -> 2799.6667
$ws.Cells.Item(80, 10).Value = 0  # J80: 4473 -> 0
$ws.Cells.Item(80, 12).Value = 0  # L80: 13419 -> 0
$ws.Cells.Item(80, 14).ClearContents()  # N80
$ws.Cells.Item(83, 8).Value = 2799.6667  # H83: 3469 -> 2799.6667
$ws.Cells.Item(83, 10).Value = 0  # J83: 4473 -> 0
$ws.Cells.Item(83, 12).Value = 0  # L83: 40257 -> 0
$ws.Cells.Item(83, 14).ClearContents()  # N83
$ws.Cells.Item(92, 8).Value = 474.33334  # H92: 549.4 -> 474.33334
$ws.Cells.Item(92, 9).Value = 269.2  # I92: 311.75 -> 269.2
$ws.Cells.Item(92, 11).Value = 807.5999999999999  # K92: 935.25 -> 807.5999999999999
$ws.Cells.Item(92, 13).Value = 440.4000000000001  # M92: 312.75 -> 440.4000000000001
$ws.Cells.Item(113, 8).Value = 522.1667  # H113: 538.82355 -> 522.1667
$ws.Cells.Item(113, 9).Value = 355  # I113: 341 -> 355
$ws.Cells.Item(113, 10).Value = 555.6  # J113: 599.6923 -> 555.6
$ws.Cells.Item(113, 11).Value = 1065  # K113: 1023 -> 1065
$ws.Cells.Item(113, 12).Value = 1666.8  # L113: 1799.0769 -> 1666.8
$ws.Cells.Item(113, 13).Value = 1105  # M113: 1147 -> 1105
$ws.Cells.Item(113, 14).Value = -6006.8  # N113: -6139.0769 -> -6006.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 0  # H99: 9950 -> 0
$ws.Cells.Item(99, 9).Value = 0  # I99: 9950 -> 0
$ws.Cells.Item(99, 11).Value = 0  # K99: 9950 -> 0
$ws.Cells.Item(99, 13).ClearContents()  # M99
$ws.Cells.Item(102, 8).Value = 3176.2942  # H102: 3127.3333 -> 3176.2942
$ws.Cells.Item(102, 9).Value = 1742.4286  # I102: 1770.2858 -> 1742.4286
$ws.Cells.Item(102, 10).Value = 4180  # J102: 3990.9092 -> 4180
$ws.Cells.Item(102, 11).Value = 1742.4286  # K102: 1770.2858 -> 1742.4286
$ws.Cells.Item(102, 12).Value = 4180  # L102: 3990.9092 -> 4180
$ws.Cells.Item(102, 13).Value = -120.4286  # M102: -148.2858000000001 -> -120.4286
$ws.Cells.Item(102, 14).Value = -7424  # N102: -7234.9092 -> -7424
$ws.Cells.Item(111, 8).Value = 42000  # H111: 0 -> 42000
$ws.Cells.Item(111, 10).Value = 42000  # J111: 0 -> 42000
$ws.Cells.Item(111, 12).Value = 42000  # L111: 0 -> 42000
$ws.Cells.Item(111, 14).Value = -48134  # N111: None -> -48134
$ws.Cells.Item(122, 8).Value = 3475.65  # H122: 3233.8333 -> 3475.65
$ws.Cells.Item(122, 9).Value = 3002.3635  # I122: 2741.6667 -> 3002.3635
$ws.Cells.Item(122, 11).Value = 9007.0905  # K122: 8225.000100000001 -> 9007.0905
$ws.Cells.Item(122, 13).Value = -6557.0905  # M122: -5775.000100000001 -> -6557.0905
$ws.Cells.Item(126, 8).Value = 10362.883  # H126: 10333.471 -> 10362.883
$ws.Cells.Item(126, 9).Value = 14416.9  # I126: 13379 -> 14416.9
$ws.Cells.Item(126, 10).Value = 4571.4287  # J126: 4750 -> 4571.4287
$ws.Cells.Item(126, 11).Value = 43250.7  # K126: 40137 -> 43250.7
$ws.Cells.Item(126, 12).Value = 13714.2861  # L126: 14250 -> 13714.2861
$ws.Cells.Item(126, 13).Value = -40780.7  # M126: -37667 -> -40780.7
$ws.Cells.Item(126, 14).Value = -18654.2861  # N126: -19190 -> -18654.2861
$ws.Cells.Item(132, 8).Value = 630031.8  # H132: 671967.3 -> 630031.8
$ws.Cells.Item(132, 9).Value = 913455.8  # I132: 1004701.5 -> 913455.8
$ws.Cells.Item(132, 11).Value = 2740367.4  # K132: 3014104.5 -> 2740367.4
$ws.Cells.Item(132, 13).Value = -2737837.4  # M132: -3011574.5 -> -2737837.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2407.5386  # H16: 2288.577 -> 2407.5386
$ws.Cells.Item(16, 9).Value = 2232  # I16: 1853 -> 2232
$ws.Cells.Item(16, 10).Value = 2612.3333  # J16: 2985.5 -> 2612.3333
$ws.Cells.Item(16, 11).Value = 2232  # K16: 1853 -> 2232
$ws.Cells.Item(16, 12).Value = 2612.3333  # L16: 2985.5 -> 2612.3333
$ws.Cells.Item(16, 13).Value = -2062  # M16: -1683 -> -2062
$ws.Cells.Item(16, 14).Value = -2952.3333  # N16: -3325.5 -> -2952.3333
$ws.Cells.Item(22, 8).Value = 10753605  # H22: 12904126 -> 10753605
$ws.Cells.Item(22, 9).Value = 1001  # I22: 0 -> 1001
$ws.Cells.Item(22, 11).Value = 1001  # K22: 0 -> 1001
$ws.Cells.Item(22, 13).Value = -706  # M22: None -> -706
$ws.Cells.Item(27, 8).Value = 10753605  # H27: 12904126 -> 10753605
$ws.Cells.Item(27, 9).Value = 1001  # I27: 0 -> 1001
$ws.Cells.Item(27, 11).Value = 1001  # K27: 0 -> 1001
$ws.Cells.Item(27, 13).Value = -894  # M27: None -> -894
$ws.Cells.Item(40, 8).Value = 3607.7144  # H40: 3425.75 -> 3607.7144
$ws.Cells.Item(40, 9).Value = 2812.875  # I40: 2680.7 -> 2812.875
$ws.Cells.Item(40, 11).Value = 2812.875  # K40: 2680.7 -> 2812.875
$ws.Cells.Item(40, 13).Value = -2676.875  # M40: -2544.7 -> -2676.875
$ws.Cells.Item(55, 8).Value = 442.21738  # H55: 428.72 -> 442.21738
$ws.Cells.Item(55, 9).Value = 452.64706  # I55: 433.78946 -> 452.64706
$ws.Cells.Item(55, 11).Value = 452.64706  # K55: 433.78946 -> 452.64706
$ws.Cells.Item(55, 13).Value = -279.64706  # M55: -260.78946 -> -279.64706
$ws.Cells.Item(122, 8).Value = 4053.3333  # H122: 4011.2 -> 4053.3333
$ws.Cells.Item(122, 9).Value = 3455  # I122: 3428.2354 -> 3455
$ws.Cells.Item(122, 11).Value = 10365  # K122: 10284.7062 -> 10365
$ws.Cells.Item(122, 13).Value = -7915  # M122: -7834.706200000001 -> -7915

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 50162.6  # H45: 54999.555 -> 50162.6
$ws.Cells.Item(45, 10).Value = 51993.125  # J45: 58473.57 -> 51993.125
$ws.Cells.Item(45, 12).Value = 51993.125  # L45: 58473.57 -> 51993.125
$ws.Cells.Item(45, 14).Value = -52975.125  # N45: -59455.57 -> -52975.125
$ws.Cells.Item(107, 8).Value = 455.0625  # H107: 404.65 -> 455.0625
$ws.Cells.Item(107, 9).Value = 444.23077  # I107: 387.47058 -> 444.23077
$ws.Cells.Item(107, 11).Value = 1332.69231  # K107: 1162.41174 -> 1332.69231
$ws.Cells.Item(107, 13).Value = 587.3076900000001  # M107: 757.58826 -> 587.3076900000001
$ws.Cells.Item(122, 8).Value = 5822.1113  # H122: 4441.643 -> 5822.1113
$ws.Cells.Item(122, 9).Value = 4666.5  # I122: 3434.818 -> 4666.5
$ws.Cells.Item(122, 11).Value = 13999.5  # K122: 10304.454 -> 13999.5
$ws.Cells.Item(122, 13).Value = -11549.5  # M122: -7854.454000000002 -> -11549.5
